# 2021-06 QLD Outbreak Paths.xlsx -- apply "Add files via upload" edit
# Adds outbreak-path rows for Q20-Q23 (Zeus Street Greek Kitchen / Household
# cluster) and refreshes the colour-gradient legend on the "Date Colours"
# sheet to match.

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet1 : outbreak path table ("Table1")
# ----------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")

# -- fix up existing rows 21 & 23 --------------------------------------
$ws1.Range("B21").Value = "Q13 w"
$ws1.Range("F21").Value = "Zeus Street Greek Kitchen"

$ws1.Range("B23").Value = "Q13 w"
$ws1.Range("F23").Value = "Zeus Street Greek Kitchen"

# -- append 4 new rows (25-28) to the table ----------------------------
$table1 = $ws1.ListObjects.Item("Table1")

$table1.ListRows.Add() | Out-Null
$ws1.Range("A25").Value = 44382
$ws1.Range("A25").NumberFormat = "d-mmm"
$ws1.Range("B25").Value = "Q13 w"
$ws1.Range("C25").Value = "Q20 w"
$ws1.Range("D25").Value = "Queensland"
$ws1.Range("F25").Value = "Carindale Greek Community Centre"
$ws1.Range("G25").Value = "Alpha (B.1.1.7)"
$ws1.Range("H25").Value = "Not isolated"

$table1.ListRows.Add() | Out-Null
$ws1.Range("A26").Value = 44382
$ws1.Range("A26").NumberFormat = "d-mmm"
$ws1.Range("B26").Value = "Q20 w"
$ws1.Range("C26").Value = "Q21 child"
$ws1.Range("D26").Value = "Queensland"
$ws1.Range("F26").Value = "Household"
$ws1.Range("G26").Value = "Alpha (B.1.1.7)"
$ws1.Range("H26").Value = "Not isolated"

$table1.ListRows.Add() | Out-Null
$ws1.Range("A27").Value = 44382
$ws1.Range("A27").NumberFormat = "d-mmm"
$ws1.Range("B27").Value = "Q17 m29"
$ws1.Range("C27").Value = "Q22 w"
$ws1.Range("D27").Value = "Queensland"
$ws1.Range("F27").Value = "Household"
$ws1.Range("G27").Value = "Alpha (B.1.1.7)"
$ws1.Range("H27").Value = "Not isolated"

$table1.ListRows.Add() | Out-Null
$ws1.Range("A28").Value = 44382
$ws1.Range("A28").NumberFormat = "d-mmm"
$ws1.Range("B28").Value = "Q17 m29"
$ws1.Range("C28").Value = "Q23"
$ws1.Range("D28").Value = "Queensland"
$ws1.Range("G28").Value = "Alpha (B.1.1.7)"
$ws1.Range("H28").Value = "Not isolated"

# ----------------------------------------------------------------------
# Sheet2 : "Date Colours" legend / gradient helper table
# ----------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Date Colours")

# -- refresh the B2:B17 gradient (shifted down a step + new final shade) --
$gradient = @("#f9f4f9","#f3e9f4","#eedfee","#e8d4e9","#e2c9e3","#dcbfdd","#d6b4d8","#d0aad2","#ca9fcd","#c495c7","#be8bc2","#b780bc","#b176b7","#ab6cb1","#a461ab")

for ($i = 0; $i -lt $gradient.Length; $i++) {
    $row = 2 + $i
    $ws2.Range("B$row").Value = $gradient[$i]
}
$ws2.Range("B17").Value = "#9e57a6"

# -- populate the new helper cells F2:U2 with the same gradient ---------
$cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws2.Range("$($cols[$i])2").Value = $gradient[$i]
}
$ws2.Range("U2").Value = "#9e57a6"

# -- selections left by the edit (Sheet1 stays the active/selected tab) --
$ws2.Range("C17").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("B23").Select() | Out-Null
